$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns before the current column E ("Folder"), shifting
# Folder, Sub-Folder, Department, Current Status, Upload Document to the right.
$ws.Range("E1:F1").EntireColumn.Insert() | Out-Null

# New header cells for the inserted columns.
$ws.Range("E1").Value = "Project Name"
$ws.Range("F1").Value = "Contract Name"

# Match header styling used by the rest of row 1 (bold Verdana 6pt font).
$ws.Range("E1:F1").Font.Bold = $true
$ws.Range("E1:F1").Font.Name = "Verdana"
$ws.Range("E1:F1").Font.Size = 6

# Column widths for the new columns (closest achievable values; the COM
# layer snaps widths to a 1/6-character grid).
$ws.Columns.Item(5).ColumnWidth = 8.35
$ws.Columns.Item(6).ColumnWidth = 8.85

# Update selection / active cell.
$ws.Range("F2").Select()
